# Quiz marksheet update: handle float/numeric marking input correctly,
# fill in the student's actual answers for question rows 16-40, and drop
# the unused duplicate "Student Ans / Correct Ans" block (columns G:H) as
# well as the now-unused second quiz block below row 18 (columns D:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12) ------------------------------------------
# Row labels now share the same bold "mtitle" style used by the header row.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# Right / Wrong / Not-Attempt / Max counts for the graded attempt.
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 28

# Marking scheme - now stored as real numbers instead of text (float input fix).
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

# Totals.
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "72/112"

# --- Student answers (rows 16-18, second quiz block D:E) ------------------
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("D18").Value = "Option B"
$ws.Range("D18").Style = "incorrectStyle"

# --- Student answers (rows 16-40, column A) --------------------------------
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"

$ws.Range("A17").Value = "Option D"
$ws.Range("A17").Style = "correctStyle"

$ws.Range("A18").Style = "normalStyle"

$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"

$ws.Range("A20").Value = "Option C"
$ws.Range("A20").Style = "incorrectStyle"

$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"

$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"

$ws.Range("A23").Style = "normalStyle"

$ws.Range("A24").Value = "Option A"
$ws.Range("A24").Style = "correctStyle"

$ws.Range("A25").Value = "Option B"
$ws.Range("A25").Style = "incorrectStyle"

$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"

$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"

$ws.Range("A28").Value = "Option D"
$ws.Range("A28").Style = "correctStyle"

$ws.Range("A29").Style = "normalStyle"

$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"

$ws.Range("A31").Value = "Option D"
$ws.Range("A31").Style = "correctStyle"

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"

$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"

$ws.Range("A34").Value = "Option B"
$ws.Range("A34").Style = "correctStyle"

$ws.Range("A35").Style = "normalStyle"

$ws.Range("A36").Value = "Option D"
$ws.Range("A36").Style = "incorrectStyle"

$ws.Range("A37").Style = "normalStyle"

$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"

$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"

$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"

# --- Drop the now-unused columns -------------------------------------------
# Second quiz block (D:E) only ever had data for rows 16-18.
$ws.Range("D19:E40").Clear()
# Duplicate third "Student Ans / Correct Ans" block (G:H) is removed entirely.
$ws.Range("G15:H40").Clear()
